# Auto-generated Excel COM-interop script
# Applies scheduled price/profit data refresh across all 8 leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 365.85715
$ws.Cells.Item(8, 9).Value = 365.85715
$ws.Cells.Item(8, 11).Value = 1097.57145
$ws.Cells.Item(8, 13).Value = -958.5714499999999
$ws.Cells.Item(11, 8).Value = 1283
$ws.Cells.Item(11, 9).Value = 1283
$ws.Cells.Item(11, 11).Value = 1283
$ws.Cells.Item(11, 13).Value = -1143
$ws.Cells.Item(17, 8).Value = 45297.523
$ws.Cells.Item(17, 10).Value = 47311.047
$ws.Cells.Item(17, 12).Value = 141933.141
$ws.Cells.Item(17, 14).Value = -142269.141
$ws.Cells.Item(21, 8).Value = 6833.3335
$ws.Cells.Item(21, 9).Value = 4200
$ws.Cells.Item(21, 10).Value = 20000
$ws.Cells.Item(21, 11).Value = 4200
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = -3732
$ws.Cells.Item(21, 14).Value = -20936
$ws.Cells.Item(23, 8).Value = 6833.3335
$ws.Cells.Item(23, 9).Value = 4200
$ws.Cells.Item(23, 10).Value = 20000
$ws.Cells.Item(23, 11).Value = 4200
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = -3966
$ws.Cells.Item(23, 14).Value = -20468
$ws.Cells.Item(51, 8).Value = 55559210
$ws.Cells.Item(51, 9).Value = 4999
$ws.Cells.Item(51, 10).Value = 83336310
$ws.Cells.Item(51, 11).Value = 4999
$ws.Cells.Item(51, 12).Value = 83336310
$ws.Cells.Item(51, 13).Value = -4515
$ws.Cells.Item(51, 14).Value = -83337278
$ws.Cells.Item(113, 8).Value = 7166.6665
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 10).Value = 9250
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 12).Value = 9250
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(113, 14).Value = -15758
$ws.Cells.Item(135, 8).Value = 975.1667
$ws.Cells.Item(135, 9).Value = 694.5454999999999
$ws.Cells.Item(135, 11).Value = 6250.9095
$ws.Cells.Item(135, 13).Value = -3715.9095
$ws.Cells.Item(141, 8).Value = 3249.2727
$ws.Cells.Item(141, 9).Value = 3454.25
$ws.Cells.Item(141, 11).Value = 10362.75
$ws.Cells.Item(141, 13).Value = -5182.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 176727.77
$ws.Cells.Item(32, 9).Value = 196481.45
$ws.Cells.Item(32, 10).Value = 8821.5
$ws.Cells.Item(32, 11).Value = 196481.45
$ws.Cells.Item(32, 12).Value = 8821.5
$ws.Cells.Item(32, 13).Value = -196194.45
$ws.Cells.Item(32, 14).Value = -9395.5
$ws.Cells.Item(38, 8).Value = 26166.666
$ws.Cells.Item(38, 10).Value = 26166.666
$ws.Cells.Item(38, 12).Value = 26166.666
$ws.Cells.Item(38, 14).Value = -27100.666
$ws.Cells.Item(74, 8).Value = 486509.22
$ws.Cells.Item(74, 9).Value = 2377.6775
$ws.Cells.Item(74, 10).Value = 1369337.4
$ws.Cells.Item(74, 11).Value = 2377.6775
$ws.Cells.Item(74, 12).Value = 1369337.4
$ws.Cells.Item(74, 13).Value = -1503.6775
$ws.Cells.Item(74, 14).Value = -1371085.4
$ws.Cells.Item(77, 8).Value = 486509.22
$ws.Cells.Item(77, 9).Value = 2377.6775
$ws.Cells.Item(77, 10).Value = 1369337.4
$ws.Cells.Item(77, 11).Value = 11888.3875
$ws.Cells.Item(77, 12).Value = 6846687
$ws.Cells.Item(77, 13).Value = -7520.387499999999
$ws.Cells.Item(77, 14).Value = -6855423

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 3050
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 5100
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 5100
$ws.Cells.Item(22, 13).Value = -827
$ws.Cells.Item(22, 14).Value = -5446
$ws.Cells.Item(99, 8).Value = 5759.55
$ws.Cells.Item(99, 9).Value = 6188.4443
$ws.Cells.Item(99, 10).Value = 1899.5
$ws.Cells.Item(99, 11).Value = 6188.4443
$ws.Cells.Item(99, 12).Value = 1899.5
$ws.Cells.Item(99, 13).Value = -4690.4443
$ws.Cells.Item(99, 14).Value = -4895.5
$ws.Cells.Item(134, 8).Value = 29034204
$ws.Cells.Item(134, 9).Value = 1691.1428
$ws.Cells.Item(134, 11).Value = 5073.428400000001
$ws.Cells.Item(134, 13).Value = -2538.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1497.0571
$ws.Cells.Item(58, 9).Value = 1344.8
$ws.Cells.Item(58, 11).Value = 1344.8
$ws.Cells.Item(58, 13).Value = -1141.8
$ws.Cells.Item(132, 8).Value = 2313.6
$ws.Cells.Item(132, 9).Value = 1683.4
$ws.Cells.Item(132, 11).Value = 5050.200000000001
$ws.Cells.Item(132, 13).Value = -2520.200000000001
$ws.Cells.Item(134, 8).Value = 2079
$ws.Cells.Item(134, 9).Value = 1798.1482
$ws.Cells.Item(134, 11).Value = 5394.444600000001
$ws.Cells.Item(134, 13).Value = -2859.444600000001
$ws.Cells.Item(136, 8).Value = 1497.0571
$ws.Cells.Item(136, 9).Value = 1344.8
$ws.Cells.Item(136, 11).Value = 4034.4
$ws.Cells.Item(136, 13).Value = -1484.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1447.7333
$ws.Cells.Item(5, 10).Value = 2500
$ws.Cells.Item(5, 12).Value = 7500
$ws.Cells.Item(5, 14).Value = -7724
$ws.Cells.Item(22, 8).Value = 19979
$ws.Cells.Item(22, 9).Value = 23723.75
$ws.Cells.Item(22, 10).Value = 5000
$ws.Cells.Item(22, 11).Value = 71171.25
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = -71002.25
$ws.Cells.Item(22, 14).Value = -15338
$ws.Cells.Item(23, 8).Value = 111357.78
$ws.Cells.Item(23, 10).Value = 143004.28
$ws.Cells.Item(23, 12).Value = 429012.84
$ws.Cells.Item(23, 14).Value = -429482.84
$ws.Cells.Item(27, 8).Value = 19979
$ws.Cells.Item(27, 9).Value = 23723.75
$ws.Cells.Item(27, 10).Value = 5000
$ws.Cells.Item(27, 11).Value = 71171.25
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = -71069.25
$ws.Cells.Item(27, 14).Value = -15204
$ws.Cells.Item(39, 8).Value = 9300
$ws.Cells.Item(39, 10).Value = 14733.333
$ws.Cells.Item(39, 12).Value = 44199.999
$ws.Cells.Item(39, 14).Value = -44787.999
$ws.Cells.Item(58, 8).Value = 5841.857
$ws.Cells.Item(58, 9).Value = 2778.6
$ws.Cells.Item(58, 11).Value = 8335.799999999999
$ws.Cells.Item(58, 13).Value = -8207.799999999999
$ws.Cells.Item(114, 8).Value = 1448.3158
$ws.Cells.Item(114, 9).Value = 1388.2858
$ws.Cells.Item(114, 11).Value = 4164.857400000001
$ws.Cells.Item(114, 13).Value = -910.8574000000008
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 1447.7333
$ws.Cells.Item(135, 10).Value = 2500
$ws.Cells.Item(135, 12).Value = 22500
$ws.Cells.Item(135, 14).Value = -27570
$ws.Cells.Item(136, 8).Value = 8221.036
$ws.Cells.Item(136, 10).Value = 15090.454
$ws.Cells.Item(136, 12).Value = 45271.362
$ws.Cells.Item(136, 14).Value = -55471.362
$ws.Cells.Item(137, 8).Value = 6766.6665
$ws.Cells.Item(137, 10).Value = 10000
$ws.Cells.Item(137, 12).Value = 30000
$ws.Cells.Item(137, 14).Value = -40200
$ws.Cells.Item(140, 8).Value = 1646.48
$ws.Cells.Item(140, 9).Value = 1485.3043
$ws.Cells.Item(140, 11).Value = 4455.9129
$ws.Cells.Item(140, 13).Value = 724.0870999999997
$ws.Cells.Item(141, 8).Value = 3125.6
$ws.Cells.Item(141, 9).Value = 3125.6
$ws.Cells.Item(141, 11).Value = 9376.799999999999
$ws.Cells.Item(141, 13).Value = -4196.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1070061.5
$ws.Cells.Item(132, 9).Value = 1226.6
$ws.Cells.Item(132, 11).Value = 3679.8
$ws.Cells.Item(132, 13).Value = -1149.8
$ws.Cells.Item(137, 8).Value = 184671.75
$ws.Cells.Item(137, 10).Value = 189596.33
$ws.Cells.Item(137, 12).Value = 189596.33
$ws.Cells.Item(137, 14).Value = -199796.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6075.1333
$ws.Cells.Item(7, 9).Value = 2239
$ws.Cells.Item(7, 10).Value = 9911.267
$ws.Cells.Item(7, 11).Value = 2239
$ws.Cells.Item(7, 12).Value = 9911.267
$ws.Cells.Item(7, 13).Value = -2127
$ws.Cells.Item(7, 14).Value = -10135.267
$ws.Cells.Item(22, 8).Value = 4584.4585
$ws.Cells.Item(22, 10).Value = 6552.1333
$ws.Cells.Item(22, 12).Value = 6552.1333
$ws.Cells.Item(22, 14).Value = -7142.1333
$ws.Cells.Item(27, 8).Value = 4584.4585
$ws.Cells.Item(27, 10).Value = 6552.1333
$ws.Cells.Item(27, 12).Value = 6552.1333
$ws.Cells.Item(27, 14).Value = -6766.1333
$ws.Cells.Item(40, 8).Value = 2706.0688
$ws.Cells.Item(40, 9).Value = 2660.3914
$ws.Cells.Item(40, 11).Value = 2660.3914
$ws.Cells.Item(40, 13).Value = -2524.3914
$ws.Cells.Item(55, 8).Value = 628.23334
$ws.Cells.Item(55, 9).Value = 358.07693
$ws.Cells.Item(55, 11).Value = 358.07693
$ws.Cells.Item(55, 13).Value = -185.07693
$ws.Cells.Item(100, 8).Value = 3734
$ws.Cells.Item(100, 9).Value = 3337.0908
$ws.Cells.Item(100, 10).Value = 3990.8235
$ws.Cells.Item(100, 11).Value = 3337.0908
$ws.Cells.Item(100, 12).Value = 3990.8235
$ws.Cells.Item(100, 13).Value = -2796.0908
$ws.Cells.Item(100, 14).Value = -5072.8235
$ws.Cells.Item(126, 8).Value = 6075.1333
$ws.Cells.Item(126, 9).Value = 2239
$ws.Cells.Item(126, 10).Value = 9911.267
$ws.Cells.Item(126, 11).Value = 6717
$ws.Cells.Item(126, 12).Value = 29733.801
$ws.Cells.Item(126, 13).Value = -4247
$ws.Cells.Item(126, 14).Value = -34673.801
$ws.Cells.Item(132, 8).Value = 3271.0667
$ws.Cells.Item(132, 9).Value = 3953.8333
$ws.Cells.Item(132, 11).Value = 11861.4999
$ws.Cells.Item(132, 13).Value = -9331.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 7145446
$ws.Cells.Item(107, 10).Value = 7145446
$ws.Cells.Item(107, 12).Value = 21436338
$ws.Cells.Item(107, 14).Value = -21440178
$ws.Cells.Item(122, 8).Value = 1329.5428
$ws.Cells.Item(122, 9).Value = 1290.25
$ws.Cells.Item(122, 10).Value = 1486.7142
$ws.Cells.Item(122, 11).Value = 3870.75
$ws.Cells.Item(122, 12).Value = 4460.142599999999
$ws.Cells.Item(122, 13).Value = -1420.75
$ws.Cells.Item(122, 14).Value = -9360.142599999999
$ws.Cells.Item(132, 8).Value = 2180.9167
$ws.Cells.Item(132, 9).Value = 1758.2
$ws.Cells.Item(132, 11).Value = 5274.6
$ws.Cells.Item(132, 13).Value = -2744.6
$ws.Cells.Item(136, 8).Value = 1227.8889
$ws.Cells.Item(136, 9).Value = 643
$ws.Cells.Item(136, 11).Value = 1929
$ws.Cells.Item(136, 13).Value = 621
